$d = $word.ActiveDocument

# Namespace fragment used for the inserted OOXML snippets below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Remember where the original content ends (right after "Detta är en
# andra ändring.") -- this is where the new, empty paragraph will need
# to be spliced in once the text paragraph below has been appended.
$originalEnd = $d.Content.End

# --- Insert the paragraph that carries the feature-branch text first,
#     appended at the (current) end of the document. Inserting a full
#     <w:p> at a Range collapsed at the very end of the story creates a
#     clean, new paragraph (as opposed to Range.InsertParagraphAfter(),
#     which leaves a stray empty <w:r> behind). ---
$endRange = $d.Content
$endRange.Collapse(0)
$textParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="PlainText"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>' +
        '</w:rPr>' +
        '<w:t>Den här texten är bara på ny-feature-branchen!</w:t>' +
    '</w:r>' +
'</w:p>'
$endRange.InsertXML($textParaXml)

# --- Now splice a second, completely empty "PlainText" paragraph (same
#     run-less shape as the existing blank paragraph already in the
#     document) in right before the paragraph just inserted, i.e. at
#     $originalEnd. Using a zero-length Range that sits strictly inside
#     the document (not at the very end of the story, which would merge
#     into the last paragraph instead of creating a new one) makes
#     InsertXML insert a sibling paragraph rather than replacing one. ---
$insertAt = $d.Range($originalEnd, $originalEnd)
$emptyParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:pStyle w:val="PlainText"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
'</w:p>'
$insertAt.InsertXML($emptyParaXml)

Write-Output "Inserted 2 new paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
